$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the raw feature codes into human-readable labels.
# (Order matters for shared-string index assignment parity with the
# source file: row 7 is written before row 6.)
$ws.Range("A4").Value = "HAS-BLED Socre"
$ws.Range("A5").Value = "High-Risk Alcohol Consumption"
$ws.Range("A7").Value = "Oral Ant. Therapy"
$ws.Range("A6").Value = "Plat. Aggr. Inhibitor Therapy"
$ws.Range("A8").Value = "Perioperative Bridging Therapy"

# Move the active selection to A13 (outside the used range), matching
# the saved cursor position in the edited workbook.
$ws.Range("A13").Select() | Out-Null
